$p = $ppt.ActivePresentation

# Delete the "Who am I?" slide (originally slide index 2).
$p.Slides.Item(2).Delete()
